$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values (which can look numeric) are written as TEXT,
# matching the original inlineStr/text representation in the sheet.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.248.64"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.858.88"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "0.7027"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "237.95"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "1.0000"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.08017"
$ws.Range("E8").Value = "  +7.84%  "
$ws.Range("D9").Value = "0.3047"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("D10").Value = "23.29"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").Value = "0.08183"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "1.867.33"
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("D13").Value = "0.7159"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "5.183"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "89.20"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "29.263.06"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "13.39"
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("D18").Value = "5.765"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").Value = "0.000007817"
$ws.Range("E19").Value = "  +2.13%  "
$ws.Range("D20").Value = "236.18"
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D22").Value = "2.106.18"
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("D23").Value = "0.9997"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "7.458"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").Value = "161.89"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "8.975"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("D27").Value = "0.1456"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "18.04"
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").Value = "1.996"
$ws.Range("E29").Value = "  +3.24%  "
$ws.Range("D30").Value = "1.435"
$ws.Range("E30").Value = "  +4.86%  "
$ws.Range("D31").Value = "1.484"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").Value = "4.403"
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("D33").Value = "4.055"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("D34").Value = "0.05217"
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("D35").Value = "1.171"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").Value = "0.7090"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Value = "2.661"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").Value = "0.01848"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").Value = "2.717"
$ws.Range("E40").Value = "  +1.76%  "
$ws.Range("D43").Value = "5.950"
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("D44").Value = "0.4280"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "70.87"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("D46").Value = "0.9994"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "103.35"
$ws.Range("E47").Value = "  +1.89%  "
$ws.Range("D48").Value = "1.783"
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("D49").Value = "2.004.69"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").Value = "9.185"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").Value = "6.945"
$ws.Range("E51").Value = "  -1.29%  "

# Row 21: only the Volume(1h) figure changed, price (D21) stayed "1.000"
$ws.Range("E21").Value = "  +0.18%  "

# Rows 41/42 swapped places (Maker <-> TrustWalletToken) with new data
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "0.9275"
$ws.Range("E41").Value = "  +3.08%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.146.95"
$ws.Range("E42").Value = "  +8.77%  "

# Restore default (General) styling on column D so no stray text-format
# style lingers on these cells (matches the original unstyled cells).
$ws.Range("D2:D51").Style = "Normal"
